$d = $word.ActiveDocument

$d.Content.Find.Execute("212×4=", $true, $false, $false, $false, $false, $true, 1, $false, "299×7=", 2) | Out-Null
$d.Content.Find.Execute("680×8=", $true, $false, $false, $false, $false, $true, 1, $false, "531×4=", 2) | Out-Null
$d.Content.Find.Execute("733×8=", $true, $false, $false, $false, $false, $true, 1, $false, "765×7=", 2) | Out-Null
$d.Content.Find.Execute("680×7=", $true, $false, $false, $false, $false, $true, 1, $false, "759×9=", 2) | Out-Null
$d.Content.Find.Execute("853×4=", $true, $false, $false, $false, $false, $true, 1, $false, "531×6=", 2) | Out-Null
$d.Content.Find.Execute("704×2=", $true, $false, $false, $false, $false, $true, 1, $false, "459×6=", 2) | Out-Null
$d.Content.Find.Execute("944×9=", $true, $false, $false, $false, $false, $true, 1, $false, "683×3=", 2) | Out-Null
$d.Content.Find.Execute("323×9=", $true, $false, $false, $false, $false, $true, 1, $false, "687×4=", 2) | Out-Null
$d.Content.Find.Execute("500×5=", $true, $false, $false, $false, $false, $true, 1, $false, "826×5=", 2) | Out-Null
$d.Content.Find.Execute("163×8=", $true, $false, $false, $false, $false, $true, 1, $false, "978×2=", 2) | Out-Null
$d.Content.Find.Execute("226×3=", $true, $false, $false, $false, $false, $true, 1, $false, "827×3=", 2) | Out-Null
$d.Content.Find.Execute("656×4=", $true, $false, $false, $false, $false, $true, 1, $false, "933×7=", 2) | Out-Null
$d.Content.Find.Execute("252×6=", $true, $false, $false, $false, $false, $true, 1, $false, "402×7=", 2) | Out-Null
$d.Content.Find.Execute("664×2=", $true, $false, $false, $false, $false, $true, 1, $false, "357×7=", 2) | Out-Null
$d.Content.Find.Execute("899×9=", $true, $false, $false, $false, $false, $true, 1, $false, "285×4=", 2) | Out-Null
$d.Content.Find.Execute("508×8=", $true, $false, $false, $false, $false, $true, 1, $false, "884×4=", 2) | Out-Null
$d.Content.Find.Execute("180×6=", $true, $false, $false, $false, $false, $true, 1, $false, "693×3=", 2) | Out-Null
$d.Content.Find.Execute("962×2=", $true, $false, $false, $false, $false, $true, 1, $false, "946×8=", 2) | Out-Null
$d.Content.Find.Execute("959×8=", $true, $false, $false, $false, $false, $true, 1, $false, "344×5=", 2) | Out-Null
$d.Content.Find.Execute("702×5=", $true, $false, $false, $false, $false, $true, 1, $false, "560×4=", 2) | Out-Null
$d.Content.Find.Execute("236×9=", $true, $false, $false, $false, $false, $true, 1, $false, "601×8=", 2) | Out-Null
$d.Content.Find.Execute("640×4=", $true, $false, $false, $false, $false, $true, 1, $false, "901×3=", 2) | Out-Null
$d.Content.Find.Execute("349×7=", $true, $false, $false, $false, $false, $true, 1, $false, "494×7=", 2) | Out-Null
$d.Content.Find.Execute("939×3=", $true, $false, $false, $false, $false, $true, 1, $false, "921×4=", 2) | Out-Null
$d.Content.Find.Execute("777×8=", $true, $false, $false, $false, $false, $true, 1, $false, "150×5=", 2) | Out-Null
